$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for the listed rows to reflect repulled data
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 4
$ws.Range("F16").Value = 6
$ws.Range("F18").Value = 2
$ws.Range("F20").Value = -2
$ws.Range("F29").Value = -5
$ws.Range("F30").Value = -4
$ws.Range("F32").Value = -2
$ws.Range("F34").Value = 5
$ws.Range("F35").Value = -3
$ws.Range("F40").Value = 0
